# Full orders excel download support added
# Populate the "заказы" worksheet with full order rows (A2:E7)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("заказы")

$data = @(
    @(1,  2, 1000,     1, 6),
    @(2,  2, 555,      1, 3),
    @(3,  2, 32313131, 1, 323),
    @(23, 1, 1000,     1, 2),
    @(24, 1, 44,       1, 2),
    @(25, 3, 400,      5, 2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Range("A$row").Value = $values[0]
    $ws.Range("B$row").Value = $values[1]
    $ws.Range("C$row").Value = $values[2]
    $ws.Range("D$row").Value = $values[3]
    $ws.Range("E$row").Value = $values[4]
}
